# Update Task tracking sheet: restructure the task table, add a new
# "Support" column, rename/merge several task descriptions & doers, and
# drop the last task row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- New "Support" column header ---------------------------------------
$ws.Range("G1").Value = "Support"

# --- Row 5 (task #4): new description, new Support doer ----------------
$ws.Range("B5").Value = "Design ui & api for account"
$ws.Range("G5").Value = "Ngân & Trí"

# --- Row 6 (task #5): new description, Doer gains a co-worker, new
#     Support doer ------------------------------------------------------
$ws.Range("B6").Value = "Design ui & api for expense & income"
$ws.Range("C6").Value = "Hạnh & Trang"
$ws.Range("G6").Value = "Vân"

# --- Row 7 (task #6): new description, new Doer, deadline cleared ------
$ws.Range("B7").Value = "Update db"
$ws.Range("C7").Value = "Trí"
$ws.Range("F7").ClearContents()

# --- Row 8 (old task #7) removed entirely; only the empty, styled
#     D:F cells remain ----------------------------------------------------
$ws.Range("A8").ClearContents()
$ws.Range("B8").ClearContents()
$ws.Range("C8").ClearContents()
$ws.Range("D8").ClearContents()

# --- Column widths: new Doer-width col C, new Support col G ------------
$ws.Columns.Item(3).ColumnWidth = 11.83
$ws.Columns.Item(7).ColumnWidth = 9.17

# --- Selection moves from the old G6 to B6 ------------------------------
[void]$ws.Range("B6").Select()
